$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Fix the header text: drop the trailing spaces in the title.
$ws.Range("A1").Value = "План на неделю с 20.02.16  по 27.02.16"

# 2) Thicken the grid border (was "medium", becomes "thick") across the
#    whole used table.
$ws.Range("A1:I8").Borders.Weight = 4

# 3) Recolor the status cells (progress highlighting):
#    D5 -> lime/teal fill
$d5 = $ws.Range("D5")
$d5.Interior.Color = 65280
$d5.Interior.PatternColor = 13421619

#    E6 -> lime/teal fill + lime green font
$e6 = $ws.Range("E6")
$e6.Interior.Color = 65280
$e6.Interior.PatternColor = 13421619
$e6.Font.Color = 65280

#    F6:H6 -> dark green/teal fill + dark green font
$f6h6 = $ws.Range("F6:H6")
$f6h6.Interior.Color = 32768
$f6h6.Interior.PatternColor = 8421376
$f6h6.Font.Color = 32768

#    H7, H8 -> red/brown fill
$hcells = $ws.Range("H7:H8")
$hcells.Interior.Color = 255
$hcells.Interior.PatternColor = 13209

# 4) Rows now carry an explicit (custom) height flag even though the
#    height values themselves are unchanged.
for ($r = 1; $r -le 8; $r++) {
  $row = $ws.Rows.Item($r)
  $row.RowHeight = $row.RowHeight
}

# 5) Selection moved from D5 to H8.
$ws.Range("H8").Select()
